$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 176 (shifts existing rows 176-267 down to 177-268)
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with the new record
$ws.Range("A176").Value = 11
$ws.Range("B176").Value = "Vega Monumental Concepción"
$ws.Range("C176").Value = "Bíobío"
$ws.Range("D176").Value = 45016
$ws.Range("E176").Value = 8
$ws.Range("F176").Value = 100112003
$ws.Range("G176").Value = "Ajo"
$ws.Range("H176").Value = "Chino"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 110
$ws.Range("K176").Value = 14000
$ws.Range("L176").Value = 15000
$ws.Range("M176").Value = 14455
$ws.Range("N176").Value = "$/caja 10 kilos"
$ws.Range("O176").Value = "China"
$ws.Range("P176").Value = 1446
$ws.Range("Q176").Value = 10
$ws.Range("R176").Value = "Hortaliza"
